$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6714.5625
$ws.Range("I18").Value = 509.2143
$ws.Range("J18").Value = 50152
$ws.Range("K18").Value = 509.2143
$ws.Range("L18").Value = 50152
$ws.Range("M18").Value = -225.2143
$ws.Range("N18").Value = -50720

# ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 8356.727999999999
$ws.Range("I21").Value = 5484.8
$ws.Range("J21").Value = 10750
$ws.Range("K21").Value = 5484.8
$ws.Range("L21").Value = 10750
$ws.Range("M21").Value = -5016.8
$ws.Range("N21").Value = -11686

# ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 8356.727999999999
$ws.Range("I23").Value = 5484.8
$ws.Range("J23").Value = 10750
$ws.Range("K23").Value = 5484.8
$ws.Range("L23").Value = 10750
$ws.Range("M23").Value = -5250.8
$ws.Range("N23").Value = -11218

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1219
$ws.Range("N29").ClearContents()

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 722975.5
$ws.Range("I58").Value = 1515498.5
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 4546495.5
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -4546345.5
$ws.Range("N58").Value = -7800

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 31251.4
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 31251.4
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 31251.4
$ws.Range("N87").Value = -33747.4

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 31251.4
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 31251.4
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 93754.20000000001
$ws.Range("N90").Value = -106234.2

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 10188.071
$ws.Range("I111").Value = 18914.5
$ws.Range("J111").Value = 3643.25
$ws.Range("K111").Value = 56743.5
$ws.Range("L111").Value = 10929.75
$ws.Range("M111").Value = -53676.5
$ws.Range("N111").Value = -17063.75

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26137.857
$ws.Range("I32").Value = 4637.25
$ws.Range("J32").Value = 88250.72
$ws.Range("K32").Value = 4637.25
$ws.Range("L32").Value = 88250.72
$ws.Range("M32").Value = -4350.25
$ws.Range("N32").Value = -88824.72

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2988.7778
$ws.Range("I61").Value = 1850
$ws.Range("J61").Value = 3314.1428
$ws.Range("K61").Value = 1850
$ws.Range("L61").Value = 3314.1428
$ws.Range("M61").Value = -1638
$ws.Range("N61").Value = -3738.1428

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4112.8857
$ws.Range("I132").Value = 4279.393
$ws.Range("J132").Value = 3446.8572
$ws.Range("K132").Value = 12838.179
$ws.Range("L132").Value = 10340.5716
$ws.Range("M132").Value = -10308.179
$ws.Range("N132").Value = -15400.5716

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2988.7778
$ws.Range("I136").Value = 1850
$ws.Range("J136").Value = 3314.1428
$ws.Range("K136").Value = 5550
$ws.Range("L136").Value = 9942.428400000001
$ws.Range("M136").Value = -3000
$ws.Range("N136").Value = -15042.4284

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4516.1816
$ws.Range("I134").Value = 5058.857
$ws.Range("J134").Value = 3566.5
$ws.Range("K134").Value = 15176.571
$ws.Range("L134").Value = 10699.5
$ws.Range("M134").Value = -12641.571
$ws.Range("N134").Value = -15769.5

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 11194.706
$ws.Range("I41").Value = 3812.5
$ws.Range("J41").Value = 13466.154
$ws.Range("K41").Value = 3812.5
$ws.Range("L41").Value = 13466.154
$ws.Range("M41").Value = -3384.5
$ws.Range("N41").Value = -14322.154

# CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 13312
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 13312
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 13312
$ws.Range("N50").Value = -14562

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7944.222
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7944.222
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7944.222
$ws.Range("N51").Value = -9416.222

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 10976
$ws.Range("I60").Value = 12000
$ws.Range("J60").Value = 10720
$ws.Range("K60").Value = 12000
$ws.Range("L60").Value = 10720
$ws.Range("M60").Value = -11489
$ws.Range("N60").Value = -11742

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 7944.222
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 7944.222
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 7944.222
$ws.Range("N61").Value = -8640.222

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1504.6666
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1504.6666
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1504.6666
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2406.6666

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 28643.5
$ws.Range("I99").Value = 4980
$ws.Range("J99").Value = 36531.332
$ws.Range("K99").Value = 4980
$ws.Range("L99").Value = 36531.332
$ws.Range("M99").Value = -3482
$ws.Range("N99").Value = -39527.332

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1780.3334
$ws.Range("I105").Value = 1353.3334
$ws.Range("J105").Value = 2207.3333
$ws.Range("K105").Value = 1353.3334
$ws.Range("L105").Value = 2207.3333
$ws.Range("M105").Value = 393.6666
$ws.Range("N105").Value = -5701.3333

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2762.1082
$ws.Range("I122").Value = 2995.1667
$ws.Range("J122").Value = 2331.8462
$ws.Range("K122").Value = 8985.500100000001
$ws.Range("L122").Value = 6995.5386
$ws.Range("M122").Value = -6535.500100000001
$ws.Range("N122").Value = -11895.5386

# CRP row 124
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 37374
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 37374
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 37374
$ws.Range("N124").Value = -42284

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 28643.5
$ws.Range("I126").Value = 4980
$ws.Range("J126").Value = 36531.332
$ws.Range("K126").Value = 14940
$ws.Range("L126").Value = 109593.996
$ws.Range("M126").Value = -12470
$ws.Range("N126").Value = -114533.996

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2000.091
$ws.Range("I34").Value = 100.666664
$ws.Range("J34").Value = 2712.375
$ws.Range("K34").Value = 301.999992
$ws.Range("L34").Value = 8137.125
$ws.Range("M34").Value = -217.999992

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1098.3334
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1098.3334
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3295.0002
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3551.0002

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11784.119
$ws.Range("I131").Value = 560
$ws.Range("J131").Value = 12708.459
$ws.Range("K131").Value = 1680
$ws.Range("L131").Value = 38125.377
$ws.Range("M131").Value = 3360
$ws.Range("N131").Value = -48205.377

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1640
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1640
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1640
$ws.Range("N113").Value = -5980

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3270744.5
$ws.Range("I126").Value = 4117.3335
$ws.Range("J126").Value = 4904058
$ws.Range("K126").Value = 12352.0005
$ws.Range("L126").Value = 14712174
$ws.Range("M126").Value = -9882.000499999998
$ws.Range("N126").Value = -14717114

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 921654.5600000001
$ws.Range("I46").Value = 1333.3334
$ws.Range("J46").Value = 1266775
$ws.Range("K46").Value = 1333.3334
$ws.Range("L46").Value = 1266775
$ws.Range("M46").Value = -1145.3334
$ws.Range("N46").Value = -1267151

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3397.1428
$ws.Range("I61").Value = 2625
$ws.Range("J61").Value = 4426.6665
$ws.Range("K61").Value = 2625
$ws.Range("L61").Value = 4426.6665
$ws.Range("M61").Value = -2423
$ws.Range("N61").Value = -4830.6665

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3397.1428
$ws.Range("I113").Value = 2625
$ws.Range("J113").Value = 4426.6665
$ws.Range("K113").Value = 2625
$ws.Range("L113").Value = 4426.6665
$ws.Range("M113").Value = -455
$ws.Range("N113").Value = -8766.666499999999

# LTW row 119
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 33798.625
$ws.Range("I119").Value = 25000
$ws.Range("J119").Value = 36731.5
$ws.Range("K119").Value = 25000
$ws.Range("L119").Value = 36731.5
$ws.Range("M119").Value = -20162
$ws.Range("N119").Value = -46407.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2999.75
$ws.Range("I136").Value = 2999.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8999.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6449.25
$ws.Range("N136").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12775.75
$ws.Range("I132").Value = 15368
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 46104
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -43574
$ws.Range("N132").Value = -20057

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1644.7667
$ws.Range("I136").Value = 455.72726
$ws.Range("J136").Value = 2333.158
$ws.Range("K136").Value = 1367.18178
$ws.Range("L136").Value = 6999.474
$ws.Range("M136").Value = 1182.81822
$ws.Range("N136").Value = -12099.474
